$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column G mirrors column F's formatting (header style, currency style, totals style)
# and width, then fills G1:G55 with the PRESUPUESTO budget values.

# Header cell G1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows G2:G54 (all zero budget values), each copying the matching F-row style
for ($r = 2; $r -le 54; $r++) {
    $ws.Range("F$r").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
    $ws.Range("G$r").Value = 0
}

# Totals row G55
$ws.Range("F55").Copy()
$ws.Range("G55").PasteSpecial(-4122)
$ws.Range("G55").Value = 0

# Column G width (raw OOXML width 17 == ColumnWidth 17 - 5/6)
$ws.Range("G1").ColumnWidth = 16.166666666666668
